$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the "Date" property value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-01-22T09:24:45+00:00"

# --- Mapping Table 1 sheet: refresh Source (A) / Target (D) columns for rows 3-15 ---
$ws = $wb.Worksheets.Item("Mapping Table 1")

$rows = @(
    @{ Row = 3;  A = "FRCDAImageIllustrative.id";                D = "FRMediaDocument.identifier" },
    @{ Row = 4;  A = "FRCDAImageIllustrative.languageCode";      D = "FRMediaDocument.content.language" },
    @{ Row = 5;  A = "FRCDAImageIllustrative.value";             D = "FRMediaDocument.content.data" },
    @{ Row = 6;  A = "FRCDAImageIllustrative.value.mediaType";   D = "FRMediaDocument.content.contentType" },
    @{ Row = 7;  A = "FRCDAImageIllustrative.subject";           D = "FRMediaDocument.subject:Patient" },
    @{ Row = 8;  A = "FRCDAImageIllustrative.specimen";          D = "FRMediaDocument.subject:Specimen" },
    @{ Row = 9;  A = "FRCDAImageIllustrative.performer";         D = "FRMediaDocument.operator.extension:performer" },
    @{ Row = 10; A = "FRCDAImageIllustrative.author";            D = "FRMediaDocument.operator.extension:author" },
    @{ Row = 11; A = "FRCDAImageIllustrative.informant";         D = "FRMediaDocument.operator.extension:informant" },
    @{ Row = 12; A = "FRCDAImageIllustrative.participant";       D = "FRMediaDocument.operator.extension:participant" },
    @{ Row = 13; A = "FRCDAImageIllustrative.entryRelationship"; D = "FRMediaDocument.basedOn" },
    @{ Row = 14; A = "FRCDAImageIllustrative.reference";         D = "FRMediaDocument.partOf" },
    @{ Row = 15; A = "FRCDAImageIllustrative.precondition";      D = "FRMediaDocument.reasonCode" }
)

foreach ($r in $rows) {
    $ws.Range("A" + $r.Row).Value = $r.A
    $ws.Range("D" + $r.Row).Value = $r.D
}
